$wb = $excel.ActiveWorkbook

# --- Rename sheet tabs ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912691014845"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912724164429"
$wb.Worksheets.Item(3).Name = "RS_TO-1650291272418403"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912724901783"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912725928879"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912690659711.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912690828063.csv"
$ws1.Range("B4").Value = "go_stims-1650291269084816.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912691014845.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16502912721017494.csv"
$ws2.Range("B3").Value = "TB-1650291272391527.csv"
$ws2.Range("B4").Value = "OB-16502912710803099.csv"
$ws2.Range("B5").Value = "ZB-match_0-165029126955121.csv"
$ws2.Range("B6").Value = "OB-16502912708630116.csv"
$ws2.Range("B7").Value = "ZB-match_4-16502912703422692.csv"
$ws2.Range("B8").Value = "ZB-match_0-16502912701698039.csv"
$ws2.Range("B9").Value = "TB-1650291271246063.csv"
$ws2.Range("B10").Value = "OB-1650291270781753.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912724431565.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912724194016.csv"
$ws4.Range("B4").Value = "MM_stims-16502912724739642.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291272444166.csv"
$ws4.Range("B6").Value = "MM_stims-16502912724901783.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912724749293.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16502912725776274.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912724947903.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502912725472598.csv"
$ws5.Range("B5").Value = "SAT_stims-16502912725206547.csv"
